# Apply fix for slopes of log transformed variables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F9:F12 values - previously "NA" text (or wrong value) now corrected numeric slopes
$ws.Range("F9").Value = 1.26728215186298
$ws.Range("F10").Value = 1.23258687060318
$ws.Range("F11").Value = 1.0962931861393299
$ws.Range("F12").Value = 0.60489880409330299

# Autofit columns B through F to match resulting "best fit" column widths
$ws.Columns("B:F").AutoFit() | Out-Null
$ws.Columns("B:B").ColumnWidth = 18
$ws.Columns("C:E").ColumnWidth = 4.833333333333333
$ws.Columns("F:F").ColumnWidth = 6.166666666666667

# Update the active selection to F12, matching the final cursor position left by the editor
$ws.Range("F12").Select() | Out-Null
